# CurtisWeeklyTaskLog_2018-06-25.docx — apply task-log / gantt update.
#
# Summary of edits:
#  1. Remove the stray "_GoBack" bookmark that sits after "25" in the
#     title paragraph.
#  2. Fill in the ID / Hours / Status cells for the (previously blank)
#     T47 row, and tighten its Status text from " Complete" to " ".
#  3. Fill in the ID / Hours / Status cells and rewrite the Description
#     for the T48 row.
#  4. Fill in the ID / Hours / Status cells for the T49 row, and rewrite
#     its Description (splitting it into two runs with a new "_GoBack"
#     bookmark placed between them, reflecting the cursor's last edit
#     position).

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function Set-EmptyCellText($cell, [string]$text) {
    # The paragraph has no runs yet (pPr only) -- assigning .Text creates
    # one, then we stamp the same sz/szCs (10pt / 20 half-points used
    # throughout this table) that every sibling run in the table uses.
    $cell.Range.Text = $text
    $cell.Range.Font.Size = 10
    $cell.Range.Font.SizeBi = 10
}

# --- 1. Drop the old "_GoBack" bookmark near "25" -------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Row for T47 (table row 3) -----------------------------------------
Set-EmptyCellText $t.Cell(3, 1) "T47"

$statusRange = $t.Cell(3, 4).Range.Duplicate
$statusRange.Find.Execute(" Complete") | Out-Null
$statusRange.Text = " "

# --- 3. Row for T48 (table row 4) -----------------------------------------
Set-EmptyCellText $t.Cell(4, 1) "T48"

$descRange = $t.Cell(4, 2).Range.Duplicate
$descRange.Find.Execute("Use Ca  ") | Out-Null
$descRange.Text = " SS           Styled all buttons"

Set-EmptyCellText $t.Cell(4, 3) "                 1.5hrs"
Set-EmptyCellText $t.Cell(4, 4) "Complete"

# --- 4. Row for T49 (table row 5) -----------------------------------------
Set-EmptyCellText $t.Cell(5, 1) "T49"

# Description cell currently holds two runs: a single space, then a tab.
# Replace the space run with the new leading text, drop a fresh
# "_GoBack" bookmark where the cursor ended up, then replace the tab run
# with the rest of the sentence.
$cell52 = $t.Cell(5, 2)
$cellStart = $cell52.Range.Start
$leadRange = $d.Range($cellStart, $cellStart + 1)
$leadRange.Text = "                Added single "

$cell52b = $t.Cell(5, 2)
$tabPos = $cell52b.Range.End - 2
$bmRange = $d.Range($tabPos, $tabPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$tailRange = $d.Range($tabPos, $tabPos + 1)
$tailRange.Text = "product page, simplified products page"

Set-EmptyCellText $t.Cell(5, 3) "3hrs"
Set-EmptyCellText $t.Cell(5, 4) "Complete"
